$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (price) ---------------------------------------------------
# Several new values look numeric ("0.9985", "304.14", ...) and a plain
# Range.Value assignment would have Excel silently convert them to numbers,
# which would both lose precision (trailing zeros) and change the stored
# cell type away from the source text cells. To avoid that (and without
# touching styles.xml / NumberFormat, which would add an unused cellXfs
# entry), each value is written as a quoted-text formula ( ="value" ) and
# then the whole range is copy/paste-special(values)-d onto itself, which
# collapses every formula down to its literal text result with no residual
# formula and no style change.
$ws.Range("D2").Formula = "=""23.610.30"""
$ws.Range("D3").Formula = "=""1.645.80"""
$ws.Range("D4").Formula = "=""0.9985"""
$ws.Range("D5").Formula = "=""0.9979"""
$ws.Range("D6").Formula = "=""304.14"""
$ws.Range("D7").Formula = "=""0.3798"""
$ws.Range("D8").Formula = "=""52.06"""
$ws.Range("D9").Formula = "=""0.3604"""
$ws.Range("D10").Formula = "=""1.243"""
$ws.Range("D11").Formula = "=""0.08190"""
$ws.Range("D12").Formula = "=""0.9970"""
$ws.Range("D14").Formula = "=""6.524"""
$ws.Range("D15").Formula = "=""7.361"""
$ws.Range("D16").Formula = "=""0.00001229"""
$ws.Range("D17").Formula = "=""1.644.82"""
$ws.Range("D18").Formula = "=""96.91"""
$ws.Range("D20").Formula = "=""6.740"""
$ws.Range("D21").Formula = "=""17.55"""
$ws.Range("D22").Formula = "=""0.9977"""
$ws.Range("D23").Formula = "=""12.56"""
$ws.Range("D24").Formula = "=""23.605.59"""
$ws.Range("D25").Formula = "=""2.524"""
$ws.Range("D26").Formula = "=""3.112"""
$ws.Range("D27").Formula = "=""21.26"""
$ws.Range("D28").Formula = "=""152.31"""
$ws.Range("D29").Formula = "=""5.187"""
$ws.Range("D30").Formula = "=""134.71"""
$ws.Range("D31").Formula = "=""1.829.37"""
$ws.Range("D32").Formula = "=""6.751"""
$ws.Range("D33").Formula = "=""1.088"""
$ws.Range("D34").Formula = "=""11.65"""
$ws.Range("D35").Formula = "=""2.044"""
$ws.Range("D36").Formula = "=""0.02793"""
$ws.Range("D37").Formula = "=""0.2511"""
$ws.Range("D38").Formula = "=""0.08812"""
$ws.Range("D40").Formula = "=""0.07018"""
$ws.Range("D41").Formula = "=""12.77"""
$ws.Range("D42").Formula = "=""0.7046"""
$ws.Range("D43").Formula = "=""1.327"""
$ws.Range("D44").Formula = "=""15.77"""
$ws.Range("D45").Formula = "=""0.6503"""
$ws.Range("D46").Formula = "=""2.334"""
$ws.Range("D47").Formula = "=""0.9974"""
$ws.Range("D48").Formula = "=""3.977"""
$ws.Range("D49").Formula = "=""0.07974"""
$ws.Range("D50").Formula = "=""127.82"""
$ws.Range("D51").Formula = "=""1.189"""
$ws.Range("D2:D51").Copy()
$ws.Range("D2:D51").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# --- Columns B/C (coin name swap on rows 34-35) --------------------------
$ws.Range("B34").Value = "FraxShare"
$ws.Range("C34").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("B35").Value = "WEMIXTOKEN"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"

# --- Column E (1h % change) ----------------------------------------------
# These strings keep their leading/trailing double spaces and a trailing
# "%" sign, which Excel never auto-parses as a number, so a direct .Value
# assignment is safe and keeps the cell as plain text.
$ws.Range("E2").Value = "  +0.75%  "
$ws.Range("E3").Value = "  +0.99%  "
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("E7").Value = "  +0.59%  "
$ws.Range("E8").Value = "  +1.01%  "
$ws.Range("E9").Value = "  -1.36%  "
$ws.Range("E10").Value = "  +1.21%  "
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("E13").Value = "  +0.13%  "
$ws.Range("E14").Value = "  -0.60%  "
$ws.Range("E15").Value = "  +0.74%  "
$ws.Range("E16").Value = "  -1.98%  "
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("E18").Value = "  +2.97%  "
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("E20").Value = "  +3.88%  "
$ws.Range("E21").Value = "  -1.01%  "
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("E23").Value = "  -1.27%  "
$ws.Range("E24").Value = "  +0.73%  "
$ws.Range("E25").Value = "  +1.87%  "
$ws.Range("E26").Value = "  -2.97%  "
$ws.Range("E27").Value = "  -0.55%  "
$ws.Range("E28").Value = "  +1.46%  "
$ws.Range("E29").Value = "  -2.13%  "
$ws.Range("E30").Value = "  +0.41%  "
$ws.Range("E31").Value = "  +1.09%  "
$ws.Range("E32").Value = "  -0.86%  "
$ws.Range("E33").Value = "  +6.99%  "
$ws.Range("E34").Value = "  +7.33%  "
$ws.Range("E35").Value = "  -10.15%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("E38").Value = "  +0.71%  "
$ws.Range("E39").Value = "  +0.48%  "
$ws.Range("E40").Value = "  -1.81%  "
$ws.Range("E41").Value = "  +4.64%  "
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("E43").Value = "  -1.89%  "
$ws.Range("E44").Value = "  -3.05%  "
$ws.Range("E45").Value = "  -0.93%  "
$ws.Range("E46").Value = "  +0.54%  "
$ws.Range("E47").Value = "  -0.27%  "
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("E49").Value = "  -0.49%  "
$ws.Range("E50").Value = "  +1.21%  "
$ws.Range("E51").Value = "  -0.86%  "
